$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number format to Text ("@") for the nombre_aides (C) and montant_total (D)
# columns on every affected row so that the new values are written as literal
# text strings (preserving formats like trailing ".00"), matching how the rest
# of the sheet stores these figures.
$targetRows = @(3,5,6,7,8,11,12,13,14,15,16,17,37,45,52,60,75,78,84,106,107,108,109,110,113,114,117,118,122,135,138,139,140,141,142,143,145,146,147,148,149,150,194,196,197,199,200,204,205,208)
foreach ($r in $targetRows) {
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(3, 3).Value = "74"
$ws.Cells.Item(3, 4).Value = "244620.00"
$ws.Cells.Item(5, 3).Value = "125"
$ws.Cells.Item(5, 4).Value = "344886.40"
$ws.Cells.Item(6, 3).Value = "348"
$ws.Cells.Item(6, 4).Value = "919010.82"
$ws.Cells.Item(7, 3).Value = "51"
$ws.Cells.Item(7, 4).Value = "125000.00"
$ws.Cells.Item(8, 3).Value = "673"
$ws.Cells.Item(8, 4).Value = "2295274.92"
$ws.Cells.Item(11, 3).Value = "32"
$ws.Cells.Item(11, 4).Value = "75270.00"
$ws.Cells.Item(12, 3).Value = "134"
$ws.Cells.Item(12, 4).Value = "375142.00"
$ws.Cells.Item(13, 3).Value = "69"
$ws.Cells.Item(13, 4).Value = "173400.00"
$ws.Cells.Item(14, 3).Value = "83"
$ws.Cells.Item(14, 4).Value = "207988.98"
$ws.Cells.Item(15, 3).Value = "17"
$ws.Cells.Item(15, 4).Value = "37593.58"
$ws.Cells.Item(16, 3).Value = "104"
$ws.Cells.Item(16, 4).Value = "423480.09"
$ws.Cells.Item(17, 3).Value = "163"
$ws.Cells.Item(17, 4).Value = "366100.00"
$ws.Cells.Item(37, 3).Value = "355"
$ws.Cells.Item(37, 4).Value = "1382857.70"
$ws.Cells.Item(45, 3).Value = "43"
$ws.Cells.Item(45, 4).Value = "138034.54"
$ws.Cells.Item(52, 3).Value = "258"
$ws.Cells.Item(52, 4).Value = "936903.67"
$ws.Cells.Item(60, 3).Value = "40"
$ws.Cells.Item(60, 4).Value = "168456.00"
$ws.Cells.Item(75, 3).Value = "39"
$ws.Cells.Item(75, 4).Value = "112579.25"
$ws.Cells.Item(78, 3).Value = "198"
$ws.Cells.Item(78, 4).Value = "541693.00"
$ws.Cells.Item(84, 3).Value = "68"
$ws.Cells.Item(84, 4).Value = "226657.55"
$ws.Cells.Item(106, 3).Value = "19"
$ws.Cells.Item(106, 4).Value = "54209.84"
$ws.Cells.Item(107, 3).Value = "67"
$ws.Cells.Item(107, 4).Value = "168310.00"
$ws.Cells.Item(108, 3).Value = "31"
$ws.Cells.Item(108, 4).Value = "91434.00"
$ws.Cells.Item(109, 3).Value = "13"
$ws.Cells.Item(109, 4).Value = "51913.61"
$ws.Cells.Item(110, 3).Value = "77"
$ws.Cells.Item(110, 4).Value = "487260.82"
$ws.Cells.Item(113, 3).Value = "22"
$ws.Cells.Item(113, 4).Value = "64267.00"
$ws.Cells.Item(114, 3).Value = "24"
$ws.Cells.Item(114, 4).Value = "59895.00"
$ws.Cells.Item(117, 3).Value = "18"
$ws.Cells.Item(117, 4).Value = "93445.92"
$ws.Cells.Item(118, 3).Value = "24"
$ws.Cells.Item(118, 4).Value = "79197.00"
$ws.Cells.Item(122, 3).Value = "238"
$ws.Cells.Item(122, 4).Value = "651508.00"
$ws.Cells.Item(135, 3).Value = "206"
$ws.Cells.Item(135, 4).Value = "574620.00"
$ws.Cells.Item(138, 3).Value = "560"
$ws.Cells.Item(138, 4).Value = "1401046.00"
$ws.Cells.Item(139, 3).Value = "1772"
$ws.Cells.Item(139, 4).Value = "4741110.93"
$ws.Cells.Item(140, 3).Value = "2379"
$ws.Cells.Item(140, 4).Value = "5950581.55"
$ws.Cells.Item(141, 3).Value = "2466"
$ws.Cells.Item(141, 4).Value = "10362367.11"
$ws.Cells.Item(142, 3).Value = "346"
$ws.Cells.Item(142, 4).Value = "976754.51"
$ws.Cells.Item(143, 3).Value = "120"
$ws.Cells.Item(143, 4).Value = "296500.00"
$ws.Cells.Item(145, 3).Value = "1007"
$ws.Cells.Item(145, 4).Value = "2614933.25"
$ws.Cells.Item(146, 3).Value = "478"
$ws.Cells.Item(146, 4).Value = "1378044.49"
$ws.Cells.Item(147, 3).Value = "364"
$ws.Cells.Item(147, 4).Value = "915700.16"
$ws.Cells.Item(148, 3).Value = "145"
$ws.Cells.Item(148, 4).Value = "357000.00"
$ws.Cells.Item(149, 3).Value = "395"
$ws.Cells.Item(149, 4).Value = "1234932.65"
$ws.Cells.Item(150, 3).Value = "835"
$ws.Cells.Item(150, 4).Value = "2011630.82"
$ws.Cells.Item(194, 3).Value = "53"
$ws.Cells.Item(194, 4).Value = "159800.00"
$ws.Cells.Item(196, 3).Value = "112"
$ws.Cells.Item(196, 4).Value = "299000.00"
$ws.Cells.Item(197, 3).Value = "349"
$ws.Cells.Item(197, 4).Value = "943788.00"
$ws.Cells.Item(199, 3).Value = "628"
$ws.Cells.Item(199, 4).Value = "2307926.16"
$ws.Cells.Item(200, 3).Value = "23"
$ws.Cells.Item(200, 4).Value = "81738.00"
$ws.Cells.Item(204, 3).Value = "62"
$ws.Cells.Item(204, 4).Value = "181926.00"
$ws.Cells.Item(205, 3).Value = "77"
$ws.Cells.Item(205, 4).Value = "181005.00"
$ws.Cells.Item(208, 3).Value = "136"
$ws.Cells.Item(208, 4).Value = "302196.77"
